$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Colors (OLE packed RGB integers)
$yellow = 65535    # RGB(255,255,0)  -> FFFF00
$orange = 49407    # RGB(255,192,0)  -> FFC000

# ---------------------------------------------------------------
# Update the text of the model-formula cells first, in the same
# order the new model descriptions were authored, so that the
# shared-string table is rebuilt in the expected sequence.
# ---------------------------------------------------------------
$ws.Range("E3").Value = "lmer(BRIX ~ treatment * year +(1|plot/plant) + (1|year:date), data = balssugboth)"
$ws.Range("E4").Value = "lmer(lnmass ~ treatment * year + (1|plot/plant) + (1|year:date), data = balssugboth)"
$ws.Range("E2").Value = "lmer(lnvol ~ treatment * year + (1|plot/plant) + (1|year:date), data = balsvolboth)"
$ws.Range("E8").Value = "lmer(BRIX ~ treatment * year + (1|plot) + (1|year:date), data = bucksugboth)"
$ws.Range("E9").Value = "lmer(lnmass ~ treatment * year + (1|plot) + (1|year:date), data = bucksugboth)"
$ws.Range("E7").Value = "lmer(lnvol ~ treatment * year + (1|plot) +(1|year:date), data = buckvolboth)"

# New "Presence/absence" section header (replaces the old row 11 footnote cell)
$ws.Range("D11").ClearContents()
$ws.Range("A11").Value = "Presence/absence"
$ws.Range("A11").Font.Bold = $true

$ws.Range("E12").Value = "glmer(necpres ~ treatment * year + (1|plot/plant) + (1| year:date), data = balsam, family = binomial)"
$ws.Range("E13").Value = "glmer(necpres ~ treatment + (1|plot) + (1|date), data = buckwt, family = binomial)"
$ws.Range("B13").Value = "NA"
$ws.Range("D13").Value = "NA"

# ---------------------------------------------------------------
# Row 2 - Balsamroot / volume
# ---------------------------------------------------------------
$ws.Range("B2").Value = 0.305
$ws.Range("C2").Value = 0.807
$ws.Range("D2").Value = 0.6016

# ---------------------------------------------------------------
# Row 3 - Balsamroot / BRIX
# ---------------------------------------------------------------
$ws.Range("B3").Value = 0.0365
$ws.Range("C3").Value = 0.3945
$ws.Range("D3").Value = 0.0932

# ---------------------------------------------------------------
# Row 4 - Balsamroot / sugar mass
# ---------------------------------------------------------------
$ws.Range("B4").Value = 0.3896
$ws.Range("C4").Interior.Color = $orange
$ws.Range("C4").Value = 0.1013
$ws.Range("D4").Interior.Color = $orange
$ws.Range("D4").Value = 0.1024

# ---------------------------------------------------------------
# Row 7 - Buckwheat / volume
# ---------------------------------------------------------------
$ws.Range("B7").Interior.Color = $yellow
$ws.Range("B7").Value = 0.0107
$ws.Range("C7").Value = 0.5593
$ws.Range("D7").Interior.Color = $orange
$ws.Range("D7").Value = 0.0771

# ---------------------------------------------------------------
# Row 8 - Buckwheat / BRIX
# ---------------------------------------------------------------
$ws.Range("B8").Value = 0.0012
$ws.Range("C8").Value = 0.2687
$ws.Range("D8").Value = 0.0181

# ---------------------------------------------------------------
# Row 9 - Buckwheat / sugar mass
# ---------------------------------------------------------------
$ws.Range("B9").Value = 0.2237
$ws.Range("C9").Value = 0.4049
$ws.Range("D9").Value = 0.2521

# ---------------------------------------------------------------
# Row 12 - Balsamroot presence/absence results
# ---------------------------------------------------------------
$ws.Range("A12").Value = "Balsamroot"
$ws.Range("B12").Interior.Color = $orange
$ws.Range("B12").Value = 0.0761
$ws.Range("C12").Value = 0.7377
$ws.Range("D12").Interior.Color = $orange
$ws.Range("D12").Value = 0.1115

# ---------------------------------------------------------------
# Row 13 - Buckwheat presence/absence results
# ---------------------------------------------------------------
$ws.Range("A13").Value = "Buckwheat"
$ws.Range("C13").Value = 0.4093

# ---------------------------------------------------------------
# Row 16 - footnote, moved down from row 11
# ---------------------------------------------------------------
$ws.Range("D16").Value = "These are LS means p-values from the indicated models"

# ---------------------------------------------------------------
# Selection, matching the saved cursor position in the workbook
# ---------------------------------------------------------------
$ws.Range("D14").Select() | Out-Null
